$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
18300.23277067938,
17912.55552228088,
17539.24431871234,
17320.42629091013,
16972.70129909401,
16631.44812727813,
17584.37173099092,
17389.78154742577,
17031.35949811108,
17273.09022877821,
16780.19201285963,
16744.54148594529,
16568.50356121331,
16334.46588526788,
16205.43440959932,
16174.44612143482,
16044.65117842302,
15927.98005195628,
15844.66941394302,
15737.81980825495,
15579.60807007752,
15426.05220739579,
15279.19274538614,
15077.42083040364,
15073.92840714972,
14971.79707693963,
14868.67298732852,
14776.63007660346,
14703.80824151496,
14626.65524763351,
14533.91394226455,
14452.58725590137,
14389.90463840106,
14387.73749508304,
14382.43484476753,
14303.62624611481,
14312.20933958654,
14240.83109155794,
14206.25703078415,
14071.46362633405,
14107.07213138274,
14074.25141069374,
14024.96077784543,
14004.74493568959,
13964.43096573496,
13855.82745079061,
13870.1347764657,
13834.23288740959,
13791.1434129584,
13731.3428529211,
13689.51329281576,
13625.98598907716,
13577.60471057701,
13519.80425924637,
13461.06049302997,
13448.26111248931,
13393.41712122999,
13349.97836176433
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $values[$i]
}

$wb.Save()